# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data updates to the Leve profit tables
# across all eight crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 185.64706
$ws.Range("I33").Value = 115.25
$ws.Range("J33").Value = 354.6
$ws.Range("K33").Value = 115.25
$ws.Range("L33").Value = 354.6
$ws.Range("M33").Value = 113.75
$ws.Range("N33").Value = -812.6

$ws.Range("H55").Value = 966.6786
$ws.Range("I55").Value = 645.5714
$ws.Range("J55").Value = 1073.7142
$ws.Range("K55").Value = 645.5714
$ws.Range("L55").Value = 1073.7142
$ws.Range("M55").Value = -431.5714
$ws.Range("N55").Value = -1501.7142

$ws.Range("H70").Value = 8496.966
$ws.Range("I70").Value = 10172
$ws.Range("J70").Value = 7136
$ws.Range("K70").Value = 30516
$ws.Range("L70").Value = 21408
$ws.Range("M70").Value = -30246
$ws.Range("N70").Value = -21948

$ws.Range("H73").Value = 8496.966
$ws.Range("I73").Value = 10172
$ws.Range("J73").Value = 7136
$ws.Range("K73").Value = 30516
$ws.Range("L73").Value = 21408
$ws.Range("M73").Value = -29580
$ws.Range("N73").Value = -23280

$ws.Range("H138").Value = 4066.8193
$ws.Range("J138").Value = 4100.721
$ws.Range("L138").Value = 12302.163
$ws.Range("N138").Value = -22582.163

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 457.10345
$ws.Range("I5").Value = 263.25
$ws.Range("J5").Value = 530.9524
$ws.Range("K5").Value = 263.25
$ws.Range("L5").Value = 530.9524
$ws.Range("M5").Value = -151.25
$ws.Range("N5").Value = -754.9524

$ws.Range("H32").Value = 10790.353
$ws.Range("I32").Value = 8564.911
$ws.Range("K32").Value = 8564.911
$ws.Range("M32").Value = -8277.911

$ws.Range("H61").Value = 1997
$ws.Range("I61").Value = 1997
$ws.Range("K61").Value = 1997
$ws.Range("M61").Value = -1785

$ws.Range("H80").Value = 21662.416
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 21662.416
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 21662.416
$ws.Range("N80").Value = -23658.416
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 21662.416
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 21662.416
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 64987.24800000001
$ws.Range("N83").Value = -74971.24800000001
$ws.Range("M83").ClearContents()

$ws.Range("H97").Value = 733.62964
$ws.Range("I97").Value = 809.4545000000001
$ws.Range("K97").Value = 809.4545000000001
$ws.Range("M97").Value = -313.4545000000001

$ws.Range("H122").Value = 2140.4644
$ws.Range("I122").Value = 2174.4814
$ws.Range("K122").Value = 6523.4442
$ws.Range("M122").Value = -4073.4442

$ws.Range("H132").Value = 7338.357
$ws.Range("I132").Value = 7895.4443
$ws.Range("K132").Value = 23686.3329
$ws.Range("M132").Value = -21156.3329

$ws.Range("H136").Value = 1997
$ws.Range("I136").Value = 1997
$ws.Range("K136").Value = 5991
$ws.Range("M136").Value = -3441

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 457.10345
$ws.Range("I4").Value = 263.25
$ws.Range("J4").Value = 530.9524
$ws.Range("K4").Value = 263.25
$ws.Range("L4").Value = 530.9524
$ws.Range("M4").Value = -148.25
$ws.Range("N4").Value = -760.9524

$ws.Range("H86").Value = 9183.166999999999
$ws.Range("I86").Value = 10700
$ws.Range("K86").Value = 10700
$ws.Range("M86").Value = -9577

$ws.Range("H89").Value = 9183.166999999999
$ws.Range("I89").Value = 10700
$ws.Range("K89").Value = 53500
$ws.Range("M89").Value = -47884

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 111111370
$ws.Range("I7").Value = 200000140
$ws.Range("J7").Value = 394
$ws.Range("K7").Value = 200000140
$ws.Range("L7").Value = 394
$ws.Range("M7").Value = -200000027
$ws.Range("N7").Value = -620

$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

$ws.Range("H51").Value = 9990
$ws.Range("I51").Value = 2475
$ws.Range("K51").Value = 2475
$ws.Range("M51").Value = -1739

$ws.Range("H61").Value = 9990
$ws.Range("I61").Value = 2475
$ws.Range("K61").Value = 2475
$ws.Range("M61").Value = -2127

$ws.Range("H132").Value = 2610.054
$ws.Range("I132").Value = 2193.9678
$ws.Range("K132").Value = 6581.903399999999
$ws.Range("M132").Value = -4051.903399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 44.55
$ws.Range("I38").Value = 69
$ws.Range("K38").Value = 207
$ws.Range("M38").Value = 140

$ws.Range("H126").Value = 15824.833
$ws.Range("I126").Value = 11666.333
$ws.Range("J126").Value = 19983.334
$ws.Range("K126").Value = 34998.999
$ws.Range("L126").Value = 59950.00199999999
$ws.Range("M126").Value = -30058.999
$ws.Range("N126").Value = -69830.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 37788.7
$ws.Range("I46").Value = 36432
$ws.Range("J46").Value = 49999
$ws.Range("K46").Value = 36432
$ws.Range("L46").Value = 49999
$ws.Range("M46").Value = -36276
$ws.Range("N46").Value = -50311

$ws.Range("H97").Value = 31292.092
$ws.Range("I97").Value = 35496.105
$ws.Range("K97").Value = 35496.105
$ws.Range("M97").Value = -35000.105

$ws.Range("H132").Value = 5825.4
$ws.Range("I132").Value = 4952.4165
$ws.Range("J132").Value = 26777
$ws.Range("K132").Value = 14857.2495
$ws.Range("L132").Value = 80331
$ws.Range("M132").Value = -12327.2495
$ws.Range("N132").Value = -85391

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 21489.957
$ws.Range("I46").Value = 32519.428
$ws.Range("J46").Value = 4333
$ws.Range("K46").Value = 32519.428
$ws.Range("L46").Value = 4333
$ws.Range("M46").Value = -32331.428
$ws.Range("N46").Value = -4709

$ws.Range("H74").Value = 101215.8
$ws.Range("I74").Value = 25215
$ws.Range("J74").Value = 120216
$ws.Range("K74").Value = 25215
$ws.Range("L74").Value = 120216
$ws.Range("M74").Value = -24217
$ws.Range("N74").Value = -122212

$ws.Range("H77").Value = 101215.8
$ws.Range("I77").Value = 25215
$ws.Range("J77").Value = 120216
$ws.Range("K77").Value = 75645
$ws.Range("L77").Value = 360648
$ws.Range("M77").Value = -70653
$ws.Range("N77").Value = -370632

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 695219.3
$ws.Range("I132").Value = 988898.8
$ws.Range("J132").Value = 4208.647
$ws.Range("K132").Value = 2966696.4
$ws.Range("L132").Value = 12625.941
$ws.Range("M132").Value = -2964166.4
$ws.Range("N132").Value = -17685.941

$ws.Range("H136").Value = 1745.5264
$ws.Range("I136").Value = 1704.5625
$ws.Range("K136").Value = 5113.6875
$ws.Range("M136").Value = -2563.6875
